$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F8").Value = 4960
$wsExhibit.Range("F11").Value = 601

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 4960
$wsAll.Range("F12").Value = 601
